$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text updates (Police Commissioner name, report volume/number, dates)
# ---------------------------------------------------------------------------
$ws.Range("M6").Value = "Jessica S. Tisch"
$ws.Range("A8").Value = "Volume 31   Number  48"
$ws.Range("C9").Value = "Report Covering the Week  11/25/2024  Through  12/1/2024"

# ---------------------------------------------------------------------------
# Crime-statistics table updates (rows 15-28)
# Cells that change from a text placeholder ("0" / "***.*") to a real numeric
# value need their NumberFormat set first so the stored style matches the
# numeric siblings in the same row.
# ---------------------------------------------------------------------------

# Row 15 - Rape
$ws.Range("M15").Value = 33.333333333333
$ws.Range("N15").Value = -38.461538461538

# Row 16 - Robbery
$ws.Range("D16").Value = 2
$ws.Range("F16").Value = 2
$ws.Range("G16").Value = 3
$ws.Range("H16").Value = -33.333333333333
$ws.Range("J16").Value = 56
$ws.Range("K16").Value = -1.785714285714
$ws.Range("L16").Value = -30.379746835443
$ws.Range("M16").Value = -29.487179487179
$ws.Range("N16").Value = -82.084690553745

# Row 17 - Fel. Assault
$ws.Range("C17").NumberFormat = "#,##0"
$ws.Range("C17").Value = 1
$ws.Range("E17").Value = -50
$ws.Range("I17").Value = 85
$ws.Range("J17").Value = 91
$ws.Range("K17").Value = -6.593406593406
$ws.Range("L17").Value = -10.526315789473
$ws.Range("M17").Value = 77.083333333333
$ws.Range("N17").Value = -32.539682539682

# Row 18 - Burglary
$ws.Range("C18").Value = 9
$ws.Range("D18").Value = 17
$ws.Range("E18").Value = -47.058823529411
$ws.Range("F18").Value = 31
$ws.Range("G18").Value = 40
$ws.Range("H18").Value = -22.5
$ws.Range("I18").Value = 286
$ws.Range("J18").Value = 300
$ws.Range("K18").Value = -4.666666666666
$ws.Range("L18").Value = 5.925925925925
$ws.Range("M18").Value = 22.746781115879
$ws.Range("N18").Value = -68.743169398907

# Row 19 - Gr. Larceny
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 12
$ws.Range("E19").Value = -58.333333333333
$ws.Range("F19").Value = 26
$ws.Range("G19").Value = 45
$ws.Range("H19").Value = -42.222222222222
$ws.Range("I19").Value = 418
$ws.Range("J19").Value = 575
$ws.Range("K19").Value = -27.304347826087
$ws.Range("L19").Value = -31.136738056013
$ws.Range("M19").Value = 22.580645161290
$ws.Range("N19").Value = -20.380952380952

# Row 20 - G.L.A.
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 6
$ws.Range("E20").Value = -83.333333333333
$ws.Range("F20").Value = 11
$ws.Range("G20").Value = 19
$ws.Range("H20").Value = -42.105263157894
$ws.Range("J20").Value = 179
$ws.Range("K20").Value = 30.167597765363
$ws.Range("L20").Value = 111.818181818182
$ws.Range("N20").Value = -92.360655737704

# Row 21 - TOTAL
$ws.Range("C21").Value = 16
$ws.Range("D21").Value = 39
$ws.Range("E21").Value = -58.974358974359
$ws.Range("F21").Value = 75
$ws.Range("G21").Value = 114
$ws.Range("H21").Value = -34.210526315789
$ws.Range("I21").Value = 1086
$ws.Range("J21").Value = 1212
$ws.Range("K21").Value = -10.396039603960
$ws.Range("L21").Value = -6.701030927835
$ws.Range("M21").Value = 30.372148859543
$ws.Range("N21").Value = -78.011743267868

# Row 24 - Petit Larceny
$ws.Range("C24").Value = 7
$ws.Range("D24").Value = 5
$ws.Range("E24").Value = 40
$ws.Range("F24").Value = 39
$ws.Range("G24").Value = 45
$ws.Range("H24").Value = -13.333333333333
$ws.Range("I24").Value = 518
$ws.Range("J24").Value = 536
$ws.Range("K24").Value = -3.358208955223
$ws.Range("L24").Value = -26.732673267326
$ws.Range("M24").Value = 24.519230769230

# Row 25 - Retail Theft
$ws.Range("C25").Value = 2
$ws.Range("D25").Value = 1
$ws.Range("E25").Value = 100
$ws.Range("F25").Value = 17
$ws.Range("G25").Value = 7
$ws.Range("H25").Value = 142.857142857143
$ws.Range("I25").Value = 121
$ws.Range("J25").Value = 112
$ws.Range("K25").Value = 8.035714285714
$ws.Range("L25").Value = 12.037037037037

# Row 26 - Misd. Assault
$ws.Range("C26").Value = 7
$ws.Range("D26").Value = 5
$ws.Range("E26").Value = 40
$ws.Range("F26").Value = 18
$ws.Range("G26").Value = 21
$ws.Range("H26").Value = -14.285714285714
$ws.Range("I26").Value = 208
$ws.Range("J26").Value = 213
$ws.Range("K26").Value = -2.347417840375
$ws.Range("L26").Value = -7.964601769911
$ws.Range("M26").Value = 30

# Row 28 - Hate Crimes
$ws.Range("D28").NumberFormat = "#,##0"
$ws.Range("D28").Value = 1
$ws.Range("E28").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E28").Value = -100
$ws.Range("F28").Value = 1
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = -50
$ws.Range("J28").Value = 17
$ws.Range("K28").Value = -29.411764705882
